# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The rule row 11 ("R40") had its Rule-name cell (B11) changed from the
# text "R40" to the literal text "1". Because "1" reads like a number,
# a leading apostrophe is used so Excel stores it as literal text (a new
# shared-string entry) instead of coercing it into a numeric cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
